# ---------------------------------------------------------------------------
# Commit: "feat: add 2022-Q4 data"
#
# 1. Insert a new worksheet named "2022-Q4" right after "总计", containing
#    the new quarter's fund-holdings table (shifts 2022-Q2 / 2022-Q1 /
#    2021-Q4 / 2021-Q3 / 2021-Q2 one position to the right).
# 2. Insert a corresponding summary row at the top of the "总计" sheet's
#    data (row 2), shifting the existing summary rows down by one and
#    renumbering the index column.
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook
$sheets = $wb.Worksheets

# --- 1. Create the new "2022-Q4" sheet, positioned right after "总计" ------

$summarySheet = $sheets.Item(1)       # "总计"
$q2Sheet      = $sheets.Item(2)       # "2022-Q2" (existing second sheet)

$newSheet = $sheets.Add()
$newSheet.Name = "2022-Q4"

# Work around a quirk where Range.Copy() into a brand-new sheet is a no-op
# unless that sheet already has at least one cell written. Seed A1, copy the
# header row + left-column styling from an existing quarter sheet (so the
# bold/border/center style matches), then clear the seed cell.
$newSheet.Range("A1").Value = "seed"
$q2Sheet.Range("B1:H1").Copy($newSheet.Range("B1:H1"))
$q2Sheet.Range("A2").Copy($newSheet.Range("A2"))
$q2Sheet.Range("A3").Copy($newSheet.Range("A3"))
$newSheet.Range("A1").ClearContents()

$newSheet.Range("B1").Value = "基金代码"
$newSheet.Range("C1").Value = "基金名称"
$newSheet.Range("D1").Value = "基金规模"
$newSheet.Range("E1").Value = "股票总仓位"
$newSheet.Range("F1").Value = "仓位占比"
$newSheet.Range("G1").Value = "持有市值(亿元)"
$newSheet.Range("H1").Value = "仓位排名"

$newSheet.Range("A2").Value = 0
$newSheet.Range("B2").Value = "'013895"
$newSheet.Range("B2").Style = "Normal"
$newSheet.Range("C2").Value = "宝盈成长精选混合A"
$newSheet.Range("D2").Value = "'8.59"
$newSheet.Range("D2").Style = "Normal"
$newSheet.Range("E2").Value = "'94.68"
$newSheet.Range("E2").Style = "Normal"
$newSheet.Range("F2").Value = "'2.90"
$newSheet.Range("F2").Style = "Normal"
$newSheet.Range("G2").Value = "'0.2491"
$newSheet.Range("G2").Style = "Normal"
$newSheet.Range("H2").Value = 9

$newSheet.Range("A3").Value = 1
$newSheet.Range("B3").Value = "'013896"
$newSheet.Range("B3").Style = "Normal"
$newSheet.Range("C3").Value = "宝盈成长精选混合C"
$newSheet.Range("D3").Value = "'3.06"
$newSheet.Range("D3").Style = "Normal"
$newSheet.Range("E3").Value = "'94.68"
$newSheet.Range("E3").Style = "Normal"
$newSheet.Range("F3").Value = "'2.90"
$newSheet.Range("F3").Style = "Normal"
$newSheet.Range("G3").Value = "'0.0887"
$newSheet.Range("G3").Style = "Normal"
$newSheet.Range("H3").Value = 9

# Move it into position 2 (right after "总计").
$newSheet.Move($sheets.Item(2))

# --- 2. Insert the new summary row into "总计" -----------------------------

$ws = $summarySheet

# Extend the table by one row, copying A6's style onto the new A7 so the
# index column keeps its formatting, then fill in the shifted values.
$ws.Range("A6").Copy($ws.Range("A7"))

$ws.Range("A7").Value = 5
$ws.Range("B7").Value = "2021-Q2"
$ws.Range("C7").Value = 2
$ws.Range("D7").Value = 1.24

$ws.Range("A6").Value = 4
$ws.Range("B6").Value = "2021-Q3"
$ws.Range("C6").Value = 9
$ws.Range("D6").Value = 2.81

$ws.Range("A5").Value = 3
$ws.Range("B5").Value = "2021-Q4"
$ws.Range("C5").Value = 7
$ws.Range("D5").Value = 1.95

$ws.Range("A4").Value = 2
$ws.Range("B4").Value = "2022-Q1"
$ws.Range("C4").Value = 2
$ws.Range("D4").Value = 0.02

$ws.Range("A3").Value = 1
$ws.Range("B3").Value = "2022-Q2"
$ws.Range("C3").Value = 2
$ws.Range("D3").Value = 0.01

$ws.Range("A2").Value = 0
$ws.Range("B2").Value = "2022-Q4"
$ws.Range("C2").Value = 2
$ws.Range("D2").Value = 0.34
